$d = $word.ActiveDocument

# 1) Insert student number after the author's name.
$rng = $d.Content
[void]$rng.Find.Execute("Krzysztof Taraszkiewicz", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" 197796")

# 2) Insert a manual line break right before "Data:" label.
$rng = $d.Content
[void]$rng.Find.Execute("Data:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)
$rng.InsertBefore([char]11)

# 3) Change the date day "13" -> "14": keep "1" as-is, change "3" to "4" (separate run).
$rng = $d.Content
[void]$rng.Find.Execute("Data: 13.12.2025", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found Data:13 Start=$($rng.Start) End=$($rng.End) Text=[$($rng.Text)]"
# Range.Start points right before 'D'; the digit '3' is at offset 7 within this match ("Data: 1" = 7 chars).
$digitRng = $d.Range($rng.Start + 7, $rng.Start + 8)
Write-Host "digitRng text=[$($digitRng.Text)]"
$digitRng.Text = "4"

Write-Host $d.Content.Text.Substring(0, 250)
